$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the Effort for the first task and compute its Remain value
$ws.Range("D2").Value = 4
$ws.Range("E2").Formula = "=C2-D2"

# Fill in the Remain formula (shared across the remaining rows) for the other tasks
$ws.Range("E3:E8").Formula = "=C3-D3"

# Update the active selection shown when the sheet is viewed
$ws.Range("C12").Select()
